$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------
$ws.Range("A2").Value = "Paul"
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0.05
$ws.Range("D2").Value = 0.05

$ws.Range("A3").Value = "Lawrence"
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.02
$ws.Range("D3").Value = 0.03

# --- Column formatting (number format + centered alignment) --------
# Columns B, C, D: percentage, centered
$ws.Columns.Item(2).HorizontalAlignment = -4108
$ws.Columns.Item(2).NumberFormat = "0%"

$ws.Columns.Item(3).HorizontalAlignment = -4108
$ws.Columns.Item(3).NumberFormat = "0%"

$ws.Columns.Item(4).HorizontalAlignment = -4108
$ws.Columns.Item(4).NumberFormat = "0%"

# Column A: text, centered
$ws.Columns.Item(1).HorizontalAlignment = -4108
$ws.Columns.Item(1).NumberFormat = "@"

# --- Column widths ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(2).ColumnWidth = 18.140625
$ws.Columns.Item(3).ColumnWidth = 9.140625
$ws.Columns.Item(4).ColumnWidth = 22.7109375

# --- Selection ---------------------------------------------------
$ws.Range("D12").Select()
